$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "ラヴァは、ハイビスカスが病気になったことを知る。姉の看病で疲れたラヴァは、目を閉じ一息ついた。目を覚ますと、それが夢だったことに気付く。病気になっていたのは、本当は自分の方だったのだ。`n"
$ws.Range("C2").Value = "Lava realizes Hibiscus is sick, and after caring for her sister, takes a moment to close her eyes and rest. When she wakes up, she finds out all of a sudden it was just a dream, and she's the one who's actually ill.`n"
